$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "CreatedAt: 2025-11-11T18:06:44"
$ws.Range("V4").Value = 91.09
$ws.Range("W4").Value = 53.73
$ws.Range("X4").Value = 80.18000000000001
$ws.Range("Y4").Value = 71.19
$ws.Range("Z4").Value = 50.08
$ws.Range("V5").Value = -103.99
$ws.Range("W5").Value = -140.05
$ws.Range("X5").Value = -9.949999999999999
$ws.Range("Y5").Value = -20.02
$ws.Range("Z5").Value = -143.58
$ws.Range("V6").Value = -15.02
$ws.Range("W6").Value = -14.34
$ws.Range("X6").Value = -5.77
$ws.Range("Y6").Value = -5.29
$ws.Range("Z6").Value = -9.1
$ws.Range("V9").Value = 89.83
$ws.Range("W9").Value = 51.76
$ws.Range("X9").Value = 84.44
$ws.Range("Y9").Value = 76.48
$ws.Range("Z9").Value = 59.79
$ws.Range("V10").Value = -103.99
$ws.Range("W10").Value = -140.05
$ws.Range("X10").Value = -9.949999999999999
$ws.Range("Y10").Value = -20.02
$ws.Range("Z10").Value = -143.58
$ws.Range("V11").Value = -16.28
$ws.Range("W11").Value = -16.3
$ws.Range("X11").Value = -1.51
$ws.Range("Z11").Value = 0.61
$ws.Range("V14").Value = 193.82
$ws.Range("W14").Value = 191.81
$ws.Range("X14").Value = 94.39
$ws.Range("Y14").Value = 96.5
$ws.Range("Z14").Value = 203.37
$ws.Range("V16").Value = -16.28
$ws.Range("W16").Value = -16.3
$ws.Range("X16").Value = -1.51
$ws.Range("Z16").Value = 0.61
$ws.Range("V19").Value = 83.11
$ws.Range("W19").Value = 49.5
$ws.Range("X19").Value = 79.76000000000001
$ws.Range("Y19").Value = 70.93000000000001
$ws.Range("Z19").Value = 49.53
$ws.Range("V20").Value = -110.35
$ws.Range("W20").Value = -142.67
$ws.Range("X20").Value = -9.949999999999999
$ws.Range("Y20").Value = -20.02
$ws.Range("Z20").Value = -143.58
$ws.Range("V21").Value = -16.64
$ws.Range("W21").Value = -15.95
$ws.Range("X21").Value = -6.19
$ws.Range("Y21").Value = -5.55
$ws.Range("Z21").Value = -9.66
$ws.Range("V24").Value = 89.47
$ws.Range("W24").Value = 52.12
$ws.Range("X24").Value = 79.76000000000001
$ws.Range("Y24").Value = 70.93000000000001
$ws.Range("Z24").Value = 49.53
$ws.Range("V25").Value = -103.99
$ws.Range("W25").Value = -140.05
$ws.Range("X25").Value = -9.949999999999999
$ws.Range("Y25").Value = -20.02
$ws.Range("Z25").Value = -143.58
$ws.Range("V26").Value = -16.64
$ws.Range("W26").Value = -15.95
$ws.Range("X26").Value = -6.19
$ws.Range("Y26").Value = -5.55
$ws.Range("Z26").Value = -9.66
$ws.Range("V29").Value = 81.17
$ws.Range("W29").Value = 47.57
$ws.Range("X29").Value = 78.93000000000001
$ws.Range("Y29").Value = 70.25
$ws.Range("Z29").Value = 48.25
$ws.Range("V30").Value = -110.35
$ws.Range("W30").Value = -142.67
$ws.Range("X30").Value = -9.949999999999999
$ws.Range("Y30").Value = -20.02
$ws.Range("Z30").Value = -143.58
$ws.Range("V31").Value = -18.58
$ws.Range("W31").Value = -17.88
$ws.Range("X31").Value = -7.02
$ws.Range("Y31").Value = -6.23
$ws.Range("Z31").Value = -10.93
$ws.Range("V34").Value = 194.54
$ws.Range("W34").Value = 125
$ws.Range("X34").Value = 97.56
$ws.Range("Y34").Value = 100.21
$ws.Range("Z34").Value = 210.33
$ws.Range("V36").Value = -15.56
$ws.Range("W36").Value = -16.66
$ws.Range("X36").Value = 1.66
$ws.Range("Y36").Value = 3.71
$ws.Range("Z36").Value = 7.57
$ws.Range("W37").Value = -66.45999999999999
$ws.Range("V39").Value = 91.09
$ws.Range("W39").Value = 53.73
$ws.Range("X39").Value = 80.18000000000001
$ws.Range("Y39").Value = 71.19
$ws.Range("Z39").Value = 50.08
$ws.Range("V40").Value = -103.99
$ws.Range("W40").Value = -140.05
$ws.Range("X40").Value = -9.949999999999999
$ws.Range("Y40").Value = -20.02
$ws.Range("Z40").Value = -143.58
$ws.Range("V41").Value = -15.02
$ws.Range("W41").Value = -14.34
$ws.Range("X41").Value = -5.77
$ws.Range("Y41").Value = -5.29
$ws.Range("Z41").Value = -9.1
$ws.Range("V44").Value = 204.18
$ws.Range("W44").Value = 203.04
$ws.Range("X44").Value = 92.93000000000001
$ws.Range("Y44").Value = 94.05
$ws.Range("Z44").Value = 200.55
$ws.Range("V46").Value = -5.92
$ws.Range("W46").Value = -5.08
$ws.Range("X46").Value = -2.97
$ws.Range("Y46").Value = -2.45
$ws.Range("V49").Value = 222.33
$ws.Range("W49").Value = 220.93
$ws.Range("X49").Value = 100.95
$ws.Range("Y49").Value = 102.77
$ws.Range("Z49").Value = 218.73
$ws.Range("V51").Value = 12.23
$ws.Range("W51").Value = 12.81
$ws.Range("X51").Value = 5.05
$ws.Range("Y51").Value = 6.27
$ws.Range("Z51").Value = 15.97
$ws.Range("V54").Value = 207.2
$ws.Range("W54").Value = 208.95
$ws.Range("X54").Value = 99.17
$ws.Range("Y54").Value = 100.94
$ws.Range("Z54").Value = 207.75
$ws.Range("V56").Value = -2.9
$ws.Range("W56").Value = 0.84
$ws.Range("X56").Value = 3.27
$ws.Range("Y56").Value = 4.44
$ws.Range("Z56").Value = 4.99
$ws.Range("W59").Value = 215.67
$ws.Range("X59").Value = 97.76000000000001
$ws.Range("Y59").Value = 99.48
$ws.Range("Z59").Value = 212.09
$ws.Range("V61").Value = 6.5
$ws.Range("W61").Value = 7.55
$ws.Range("X61").Value = 1.86
$ws.Range("Y61").Value = 2.98
$ws.Range("Z61").Value = 9.33
$ws.Range("V64").Value = 220.69
$ws.Range("W64").Value = 220.47
$ws.Range("X64").Value = 99.69
$ws.Range("Y64").Value = 101.47
$ws.Range("Z64").Value = 215.7
$ws.Range("V66").Value = 10.59
$ws.Range("W66").Value = 12.35
$ws.Range("X66").Value = 3.79
$ws.Range("Y66").Value = 4.97
$ws.Range("Z66").Value = 12.94
$ws.Range("W69").Value = 220
$ws.Range("X69").Value = 98.66
$ws.Range("Y69").Value = 101.26
$ws.Range("Z69").Value = 217.09
$ws.Range("V71").Value = 9.9
$ws.Range("W71").Value = 11.88
$ws.Range("X71").Value = 2.76
$ws.Range("Y71").Value = 4.76
$ws.Range("Z71").Value = 14.33
$ws.Range("W74").Value = 218.15
$ws.Range("X74").Value = 98.87
$ws.Range("Y74").Value = 100.52
$ws.Range("Z74").Value = 213.88
$ws.Range("V76").Value = 8.529999999999999
$ws.Range("W76").Value = 10.04
$ws.Range("X76").Value = 2.97
$ws.Range("Y76").Value = 4.02
$ws.Range("Z76").Value = 11.12
$ws.Range("V79").Value = 219.29
$ws.Range("W79").Value = 219.31
$ws.Range("X79").Value = 99.39
$ws.Range("Y79").Value = 101.12
$ws.Range("Z79").Value = 214.79
$ws.Range("V81").Value = 9.19
$ws.Range("W81").Value = 11.19
$ws.Range("X81").Value = 3.49
$ws.Range("Y81").Value = 4.62
$ws.Range("Z81").Value = 12.03
$ws.Range("V84").Value = 200.86
$ws.Range("W84").Value = 206.26
$ws.Range("X84").Value = 98.06
$ws.Range("Y84").Value = 99.59
$ws.Range("Z84").Value = 198.78
$ws.Range("V86").Value = -9.24
$ws.Range("W86").Value = -1.86
$ws.Range("X86").Value = 2.16
$ws.Range("Y86").Value = 3.09
$ws.Range("Z86").Value = -3.98
$ws.Range("V89").Value = 87.53
$ws.Range("W89").Value = 50.19
$ws.Range("X89").Value = 78.93000000000001
$ws.Range("Y89").Value = 70.25
$ws.Range("Z89").Value = 48.25
$ws.Range("V90").Value = -103.99
$ws.Range("W90").Value = -140.05
$ws.Range("X90").Value = -9.949999999999999
$ws.Range("Y90").Value = -20.02
$ws.Range("Z90").Value = -143.58
$ws.Range("V91").Value = -18.58
$ws.Range("W91").Value = -17.88
$ws.Range("X91").Value = -7.02
$ws.Range("Y91").Value = -6.23
$ws.Range("Z91").Value = -10.93
